$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 85, pushing the existing rows 85-98
# down to 87-100 (dimension grows from A1:R98 to A1:R100).
$ws.Rows.Item(85).Insert()
$ws.Rows.Item(85).Insert()

# --- New row 85 ---
$ws.Cells.Item(85,1).Value  = 1
$ws.Cells.Item(85,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(85,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(85,4).Value  = 45154
$ws.Cells.Item(85,5).Value  = 15
$ws.Cells.Item(85,6).Value  = 100112045
$ws.Cells.Item(85,7).Value  = "Zapallo"
$ws.Cells.Item(85,8).Value  = "Camote"
$ws.Cells.Item(85,9).Value  = "1a nueva(o)"
$ws.Cells.Item(85,10).Value = 700
$ws.Cells.Item(85,11).Value = 600
$ws.Cells.Item(85,12).Value = 630
$ws.Cells.Item(85,13).Value = 615
$ws.Cells.Item(85,14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(85,15).Value = "Perú"
$ws.Cells.Item(85,16).Value = 615
$ws.Cells.Item(85,17).Value = 1
$ws.Cells.Item(85,18).Value = "Hortaliza"

# --- New row 86 ---
$ws.Cells.Item(86,1).Value  = 1
$ws.Cells.Item(86,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(86,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(86,4).Value  = 45154
$ws.Cells.Item(86,5).Value  = 15
$ws.Cells.Item(86,6).Value  = 100112045
$ws.Cells.Item(86,7).Value  = "Zapallo"
$ws.Cells.Item(86,8).Value  = "Camote"
$ws.Cells.Item(86,9).Value  = "2a nueva(o)"
$ws.Cells.Item(86,10).Value = 900
$ws.Cells.Item(86,11).Value = 550
$ws.Cells.Item(86,12).Value = 580
$ws.Cells.Item(86,13).Value = 565
$ws.Cells.Item(86,14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(86,15).Value = "Perú"
$ws.Cells.Item(86,16).Value = 565
$ws.Cells.Item(86,17).Value = 1
$ws.Cells.Item(86,18).Value = "Hortaliza"
